$d = $word.ActiveDocument

# 1) First paragraph currently reads "1S/256MB" with the _GoBack bookmark.
#    It must become two runs: "2" (eastAsia-hinted) + "S/256MB", with the
#    _GoBack bookmark removed from here.
$p1 = $d.Paragraphs(1).Range

# Replace the leading "1" with "2"
$rng = $p1.Duplicate
$rng.Find.Execute("1S/256MB", $true, $false, $false, $false, $false, $true, 1, $false, "2S/256MB", 2)

# Now split "2S/256MB" into two runs: "2" and "S/256MB", giving "2" the
# eastAsia rFonts hint like the "D" run in the next paragraph.
$para1 = $d.Paragraphs(1)
$firstRun = $para1.Range.Characters(1)
$firstRun.Font.NameFarEast = $firstRun.Font.NameFarEast

# Remove the _GoBack bookmark from paragraph 1 if present
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Add the _GoBack bookmark at the end of paragraph 2 ("Description")
$para2 = $d.Paragraphs(2)
$end2 = $para2.Range.End - 1
$bmRange = $d.Range($end2, $end2)
$d.Bookmarks.Add("_GoBack", $bmRange)
